$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

function Rename-InlineShape($range, $newName) {
    $ishp = $range.InlineShapes.Item(1)
    $shp = $ishp.ConvertToShape()
    $shp.Name = $newName
    [void]$shp.ConvertToInlineShape()
}

# Footer (default, type=2/footer2.xml): PearsonLogo.png  image1.png -> image2.png
Rename-InlineShape $sec.Footers.Item(1).Range "image2.png"

# Footer (first page, type=1/footer1.xml): PearsonLogo.png  image1.png -> image2.png
Rename-InlineShape $sec.Footers.Item(2).Range "image2.png"

# Header (first page, type=1/header1.xml): BTec_Logo-Orange  image2.jpg -> image1.jpg
Rename-InlineShape $sec.Headers.Item(2).Range "image1.jpg"
